$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1907.5397
$ws.Range("I132").Value = 1793.037
$ws.Range("J132").Value = 2594.5557
$ws.Range("K132").Value = 5379.111
$ws.Range("L132").Value = 7783.6671
$ws.Range("M132").Value = -2849.111
$ws.Range("N132").Value = -12843.6671

$ws.Range("H137").Value = 715.8929
$ws.Range("I137").Value = 663.2692
$ws.Range("J137").Value = 1400
$ws.Range("K137").Value = 1989.8076
$ws.Range("L137").Value = 4200
$ws.Range("M137").Value = 560.1924000000001
$ws.Range("N137").Value = -9300

$ws.Range("H138").Value = 2970
$ws.Range("I138").Value = 649.58826
$ws.Range("J138").Value = 4165.364
$ws.Range("K138").Value = 1948.76478
$ws.Range("L138").Value = 12496.092
$ws.Range("M138").Value = 3191.23522
$ws.Range("N138").Value = -22776.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1160
$ws.Range("I61").Value = 950
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 950
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -738
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 1044.875
$ws.Range("I74").Value = 1127.1052
$ws.Range("J74").Value = 924.6923
$ws.Range("K74").Value = 1127.1052
$ws.Range("L74").Value = 924.6923
$ws.Range("M74").Value = -253.1052
$ws.Range("N74").Value = -2672.6923

$ws.Range("H77").Value = 1044.875
$ws.Range("I77").Value = 1127.1052
$ws.Range("J77").Value = 924.6923
$ws.Range("K77").Value = 5635.526
$ws.Range("L77").Value = 4623.4615
$ws.Range("M77").Value = -1267.526
$ws.Range("N77").Value = -13359.4615

$ws.Range("H132").Value = 2242.5293
$ws.Range("I132").Value = 858.8571
$ws.Range("J132").Value = 3211.1
$ws.Range("K132").Value = 2576.5713
$ws.Range("L132").Value = 9633.3
$ws.Range("M132").Value = -46.57129999999961
$ws.Range("N132").Value = -14693.3

$ws.Range("H136").Value = 1160
$ws.Range("I136").Value = 950
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2850
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -300
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 90041.95
$ws.Range("I134").Value = 2855.3125
$ws.Range("J134").Value = 289325.72
$ws.Range("K134").Value = 8565.9375
$ws.Range("L134").Value = 867977.1599999999
$ws.Range("M134").Value = -6030.9375
$ws.Range("N134").Value = -873047.1599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2157.6072
$ws.Range("I31").Value = 2075.4
$ws.Range("J31").Value = 2842.6667
$ws.Range("K31").Value = 2075.4
$ws.Range("L31").Value = 2842.6667
$ws.Range("M31").Value = -1780.4
$ws.Range("N31").Value = -3432.6667

$ws.Range("H34").Value = 2157.6072
$ws.Range("I34").Value = 2075.4
$ws.Range("J34").Value = 2842.6667
$ws.Range("K34").Value = 2075.4
$ws.Range("L34").Value = 2842.6667
$ws.Range("M34").Value = -1873.4
$ws.Range("N34").Value = -3246.6667

$ws.Range("H58").Value = 3640.4443
$ws.Range("I58").Value = 816.8182
$ws.Range("J58").Value = 8077.5713
$ws.Range("K58").Value = 816.8182
$ws.Range("L58").Value = 8077.5713
$ws.Range("M58").Value = -613.8182
$ws.Range("N58").Value = -8483.5713

$ws.Range("H132").Value = 1489.8334
$ws.Range("I132").Value = 920.7895
$ws.Range("J132").Value = 2472.7273
$ws.Range("K132").Value = 2762.3685
$ws.Range("L132").Value = 7418.1819
$ws.Range("M132").Value = -232.3685
$ws.Range("N132").Value = -12478.1819

$ws.Range("H134").Value = 2088.175
$ws.Range("I134").Value = 1566.4688
$ws.Range("J134").Value = 4175
$ws.Range("K134").Value = 4699.4064
$ws.Range("L134").Value = 12525
$ws.Range("M134").Value = -2164.4064
$ws.Range("N134").Value = -17595

$ws.Range("H136").Value = 3640.4443
$ws.Range("I136").Value = 816.8182
$ws.Range("J136").Value = 8077.5713
$ws.Range("K136").Value = 2450.4546
$ws.Range("L136").Value = 24232.7139
$ws.Range("M136").Value = 99.54539999999997
$ws.Range("N136").Value = -29332.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4902098
$ws.Range("I2").Value = 11764723
$ws.Range("J2").Value = 223.14285
$ws.Range("K2").Value = 70588338
$ws.Range("L2").Value = 1338.8571
$ws.Range("M2").Value = -70588225
$ws.Range("N2").Value = -1564.8571

$ws.Range("H34").Value = 3301.4285
$ws.Range("I34").Value = 518.3333
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 1554.9999
$ws.Range("L34").Value = 60000
$ws.Range("M34").Value = -1470.9999
$ws.Range("N34").Value = -60168

$ws.Range("H39").Value = 2018.55
$ws.Range("J39").Value = 2018.55
$ws.Range("L39").Value = 6055.65
$ws.Range("N39").Value = -6643.65

$ws.Range("H55").Value = 48558.855
$ws.Range("I55").Value = 250524.75
$ws.Range("J55").Value = 1037.4706
$ws.Range("K55").Value = 751574.25
$ws.Range("L55").Value = 3112.4118
$ws.Range("M55").Value = -751397.25
$ws.Range("N55").Value = -3466.4118

$ws.Range("H92").Value = 312.25
$ws.Range("J92").Value = 312.25
$ws.Range("L92").Value = 936.75
$ws.Range("N92").Value = -3432.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5150
$ws.Range("I70").Value = 4646.154
$ws.Range("J70").Value = 7333.3335
$ws.Range("K70").Value = 4646.154
$ws.Range("L70").Value = 7333.3335
$ws.Range("M70").Value = -4376.154
$ws.Range("N70").Value = -7873.3335

$ws.Range("H73").Value = 5150
$ws.Range("I73").Value = 4646.154
$ws.Range("J73").Value = 7333.3335
$ws.Range("K73").Value = 4646.154
$ws.Range("L73").Value = 7333.3335
$ws.Range("M73").Value = -3710.154
$ws.Range("N73").Value = -9205.3335

$ws.Range("H132").Value = 3807.6843
$ws.Range("I132").Value = 4441.75
$ws.Range("J132").Value = 3346.5454
$ws.Range("K132").Value = 13325.25
$ws.Range("L132").Value = 10039.6362
$ws.Range("M132").Value = -10795.25
$ws.Range("N132").Value = -15099.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2469.7273
$ws.Range("I132").Value = 1402.6923
$ws.Range("J132").Value = 4011
$ws.Range("K132").Value = 4208.0769
$ws.Range("L132").Value = 12033
$ws.Range("M132").Value = -1678.0769
$ws.Range("N132").Value = -17093

$ws.Range("H136").Value = 3836.9355
$ws.Range("I136").Value = 1309.4231
$ws.Range("J136").Value = 16980
$ws.Range("K136").Value = 3928.2693
$ws.Range("L136").Value = 50940
$ws.Range("M136").Value = -1378.2693
$ws.Range("N136").Value = -56040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5330.8335
$ws.Range("J45").Value = 5761.75
$ws.Range("L45").Value = 5761.75
$ws.Range("N45").Value = -6743.75

$ws.Range("H132").Value = 2077.2666
$ws.Range("I132").Value = 1856.2
$ws.Range("J132").Value = 2519.4
$ws.Range("K132").Value = 5568.6
$ws.Range("L132").Value = 7558.200000000001
$ws.Range("M132").Value = -3038.6
$ws.Range("N132").Value = -12618.2

$ws.Range("H136").Value = 2226.182
$ws.Range("I136").Value = 2225.3
$ws.Range("J136").Value = 2235
$ws.Range("K136").Value = 6675.900000000001
$ws.Range("L136").Value = 6705
$ws.Range("M136").Value = -4125.900000000001
$ws.Range("N136").Value = -11805
